$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to remain text even when the value looks like a
    # plain number (e.g. "593.95" or "1.00"), mirroring the source data
    # where these price strings are stored as text, not numbers.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value2 = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value2 = "68.675.00"
$ws.Range("E2").Value2 = "  +2.38%  "

# Row 3 - Ethereum
$ws.Range("D3").Value2 = "2.533.23"
$ws.Range("E3").Value2 = "  +2.69%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value2 = "  +0.03%  "

# Row 5 - BNB
Set-TextValue "D5" "593.95"
$ws.Range("E5").Value2 = "  +1.94%  "

# Row 6 - Solana
Set-TextValue "D6" "178.00"
$ws.Range("E6").Value2 = "  +2.35%  "

# Row 8 - XRP
Set-TextValue "D8" "0.518"
$ws.Range("E8").Value2 = "  +1.24%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value2 = "2.532.34"
$ws.Range("E9").Value2 = "  +2.69%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.145"
$ws.Range("E10").Value2 = "  +6.20%  "

# Row 11 - TRON
$ws.Range("E11").Value2 = "  -0.97%  "

# Row 12 - Toncoin
$ws.Range("E12").Value2 = "  +1.26%  "

# Row 13 - Cardano
$ws.Range("E13").Value2 = "  +2.05%  "

# Row 15 - Avalanche
Set-TextValue "D15" "26.14"
$ws.Range("E15").Value2 = "  +3.14%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value2 = "68.308.95"
$ws.Range("E16").Value2 = "  +2.07%  "

# Row 17 - ShibaInu
Set-TextValue "D17" "0.0000171"
$ws.Range("E17").Value2 = "  +1.53%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value2 = "2.519.97"
$ws.Range("E18").Value2 = "  +2.13%  "

# Row 19 - Chainlink
Set-TextValue "D19" "11.12"
$ws.Range("E19").Value2 = "  +2.27%  "

# Row 20 - Uniswap
Set-TextValue "D20" "7.52"
$ws.Range("E20").Value2 = "  +1.07%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "353.70"
$ws.Range("E21").Value2 = "  +1.67%  "

# Row 22 - Polkadot
$ws.Range("E22").Value2 = "  +5.03%  "

# Row 23 - was Litecoin, now Dai
$ws.Range("B23").Value2 = "Dai"
$ws.Range("C23").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D23" "1.00"
$ws.Range("E23").Value2 = "  -0.02%  "

# Row 24 - was Dai, now Litecoin
$ws.Range("B24").Value2 = "Litecoin"
$ws.Range("C24").Value2 = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D24" "70.98"
$ws.Range("E24").Value2 = "  +2.40%  "

# Row 25 - NEARProtocol
$ws.Range("E25").Value2 = "  +0.94%  "

# Row 26 - SuiNetwork
$ws.Range("E26").Value2 = "  -4.13%  "

# Row 27 - Aptos
$ws.Range("E27").Value2 = "  -1.02%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value2 = "2.652.05"
$ws.Range("E28").Value2 = "  +2.28%  "

# Row 29 - Binance-PegBSC-USD
Set-TextValue "D29" "0.999"
$ws.Range("E29").Value2 = "  +0.03%  "

# Row 30 - Bittensor
Set-TextValue "D30" "514.52"
$ws.Range("E30").Value2 = "  +3.40%  "

# Row 31 - PEPE
$sub3 = [string][char]0x2083
$ws.Range("D31").Value2 = "0.0" + $sub3 + "0900"
$ws.Range("E31").Value2 = "  +0.37%  "

# Row 32 - InternetComputer(DFINITY)
Set-TextValue "D32" "7.82"
$ws.Range("E32").Value2 = "  +1.29%  "

# Row 33 - Fetch.AI
$ws.Range("E33").Value2 = "  +2.74%  "

# Row 34 - PancakeSwap
$ws.Range("E34").Value2 = "  +1.54%  "

# Row 37 - Kaspa
$ws.Range("E37").Value2 = "  +0.32%  "

# Row 38 - EthereumClassic
Set-TextValue "D38" "18.44"
$ws.Range("E38").Value2 = "  +1.90%  "

# Row 39 - WhiteBITCoin
Set-TextValue "D39" "18.69"
$ws.Range("E39").Value2 = "  +0.09%  "

# Row 40 - ImmutableX
$ws.Range("E40").Value2 = "  +0.38%  "

# Row 41 - Stacks
Set-TextValue "D41" "1.76"
$ws.Range("E41").Value2 = "  +4.86%  "

# Row 42 - USDe
$ws.Range("E42").Value2 = "  +0.03%  "

# Row 43 - RenderToken
Set-TextValue "D43" "4.86"
$ws.Range("E43").Value2 = "  +1.19%  "

# Row 44 - PolygonEcosystemToken
$ws.Range("E44").Value2 = "  +0.48%  "

# Row 45 - dogwifhat
$ws.Range("E45").Value2 = "  +2.34%  "

# Row 46 - Aave
Set-TextValue "D46" "152.68"
$ws.Range("E46").Value2 = "  +7.31%  "

# Row 47 - Filecoin
$ws.Range("E47").Value2 = "  +3.31%  "

# Row 48 - BabyDogeCoin
$ws.Range("E48").Value2 = "  +3.57%  "

# Row 49 - ARBITRUM
Set-TextValue "D49" "0.521"
$ws.Range("E49").Value2 = "  +2.68%  "

# Row 50 - Optimism
$ws.Range("E50").Value2 = "  +3.58%  "

# Row 51 - Cronos
$ws.Range("E51").Value2 = "  +0.59%  "
